$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column holds text-formatted numbers (e.g. thousand-dot
# separated strings like "27.811.77", or values with significant
# trailing zeros like "41.50"). Force the Text format first so
# Excel does not auto-convert these into numeric values and mangle
# them (dropping trailing zeros, re-parsing multi-dot strings, etc).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.811.77'
$ws.Range("E2").Value = '  -3.17%  '
$ws.Range("D3").Value = '1.796.41'
$ws.Range("E3").Value = '  -0.63%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '315.57'
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").Value = '0.5360'
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '0.3847'
$ws.Range("E8").Value = '  +1.45%  '
$ws.Range("D9").Value = '0.07448'
$ws.Range("E9").Value = '  -1.15%  '
$ws.Range("D10").Value = '41.50'
$ws.Range("E10").Value = '  -2.65%  '
$ws.Range("D11").Value = '1.086'
$ws.Range("E11").Value = '  -2.95%  '
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  -0.06%  '
$ws.Range("D13").Value = '6.249'
$ws.Range("E13").Value = '  +0.90%  '
$ws.Range("D14").Value = '7.468'
$ws.Range("E14").Value = '  +1.04%  '
$ws.Range("D15").Value = '20.37'
$ws.Range("E15").Value = '  -2.77%  '
$ws.Range("D16").Value = '1.790.14'
$ws.Range("E16").Value = '  -0.81%  '
$ws.Range("D17").Value = '88.43'
$ws.Range("E17").Value = '  -2.46%  '
$ws.Range("D18").Value = '0.00001060'
$ws.Range("E18").Value = '  -0.67%  '
$ws.Range("D19").Value = '0.06525'
$ws.Range("E19").Value = '  +1.02%  '
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.06%  '
$ws.Range("E21").Value = '  +0.37%  '
$ws.Range("D22").Value = '5.972'
$ws.Range("E22").Value = '  +0.92%  '
$ws.Range("D23").Value = '27.842.38'
$ws.Range("E23").Value = '  -3.09%  '
$ws.Range("D24").Value = '11.11'
$ws.Range("E24").Value = '  -0.88%  '
$ws.Range("D25").Value = '2.091'
$ws.Range("E25").Value = '  -0.72%  '
$ws.Range("D26").Value = '156.64'
$ws.Range("E26").Value = '  -2.70%  '
$ws.Range("E27").Value = '  -1.00%  '
$ws.Range("D28").Value = '2.000.60'
$ws.Range("E28").Value = '  -0.67%  '
$ws.Range("D29").Value = '2.335'
$ws.Range("E29").Value = '  -1.22%  '
$ws.Range("D30").Value = '121.67'
$ws.Range("E30").Value = '  -1.19%  '
$ws.Range("E31").Value = '  +0.98%  '
$ws.Range("D32").Value = '0.1093'
$ws.Range("E32").Value = '  +3.00%  '
$ws.Range("D33").Value = '3.654'
$ws.Range("E33").Value = '  -1.00%  '
$ws.Range("D34").Value = '5.528'
$ws.Range("E34").Value = '  -2.45%  '
$ws.Range("D35").Value = '0.07015'
$ws.Range("E35").Value = '  +8.78%  '
$ws.Range("D36").Value = '0.2199'
$ws.Range("E36").Value = '  -2.82%  '
$ws.Range("D37").Value = '0.02276'
$ws.Range("E37").Value = '  -1.73%  '
$ws.Range("D38").Value = '5.086'
$ws.Range("E38").Value = '  +0.65%  '
$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").Value = '11.40'
$ws.Range("E39").Value = '  +0.88%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '8.445'
$ws.Range("E40").Value = '  -3.76%  '
$ws.Range("D41").Value = '0.6114'
$ws.Range("E41").Value = '  -2.29%  '
$ws.Range("B42").Value = 'WEMIXTOKEN'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D42").Value = '1.412'
$ws.Range("E42").Value = '  +0.96%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '1.162'
$ws.Range("E43").Value = '  -4.37%  '
$ws.Range("D44").Value = '13.26'
$ws.Range("E44").Value = '  -0.33%  '
$ws.Range("D45").Value = '3.677'
$ws.Range("E45").Value = '  -0.34%  '
$ws.Range("D46").Value = '0.5717'
$ws.Range("E46").Value = '  -2.82%  '
$ws.Range("D47").Value = '125.05'
$ws.Range("E47").Value = '  -0.93%  '
$ws.Range("D48").Value = '1.911'
$ws.Range("E48").Value = '  -2.36%  '
$ws.Range("E49").Value = '  +1.71%  '
$ws.Range("D50").Value = '0.06789'
$ws.Range("E50").Value = '  -1.46%  '
$ws.Range("D51").Value = '71.81'
$ws.Range("E51").Value = '  -1.40%  '
